# Update "存货增加" sheet: drop the oldest four year rows (2000/2002/2005/2007)
# and append a new 2020 row after the remaining 2010/2012/2015/2017 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 2-5 (2000年, 2002年, 2005年, 2007年); remaining rows shift up,
# so old row 6 (2010年) becomes new row 2, ..., old row 9 (2017年) becomes new row 5.
$ws.Range("A2:A5").EntireRow.Delete()

# Copy formatting from row 2 (label cell A2, which carries the bold/border/center
# style) onto the new row 6 label cell before filling in the 2020 data.
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("A6").Value = "2020年"
$ws.Range("B6").Value = -2103816.76978438
$ws.Range("D6").Value = 54901501.621725
$ws.Range("F6").Value = 91509937.67722151
$ws.Range("J6").Value = 37249209.8747197
$ws.Range("K6").Value = -25605149.2034971
$ws.Range("L6").Value = 1785150.09795483
$ws.Range("M6").Value = -3571111.18213603
$ws.Range("O6").Value = 462218.578754962
$ws.Range("P6").Value = 12957019.0619067
$ws.Range("R6").Value = 2140696.25007459
$ws.Range("S6").Value = 8525640.232746361
